# Append the new match row (row 60) to the Azerbaijan Premier League sheet,
# mirroring the existing table layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 60
$prevRow = 59

# Copy the formatting of the previous data row onto the new row first so the
# appended cells inherit the same number formats / styles as the rest of the
# table (bold+bordered index column, datetime format on the match-date
# column, etc).
$ws.Range("A$prevRow`:V$prevRow").Copy() | Out-Null
$ws.Range("A$newRow`:V$newRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 59
$ws.Cells.Item($newRow, 2).Value = "azerbaijan"
$ws.Cells.Item($newRow, 3).Value = "premier-league"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45235.64583333334
$ws.Cells.Item($newRow, 6).Value = "Sabah Baku"
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = "Zira"
$ws.Cells.Item($newRow, 9).Value = 1
$ws.Cells.Item($newRow, 10).Value = 1.57
$ws.Cells.Item($newRow, 11).Value = "04/11/2023 03:43"
$ws.Cells.Item($newRow, 12).Value = 1.9
$ws.Cells.Item($newRow, 13).Value = "05/11/2023 15:17"
$ws.Cells.Item($newRow, 14).Value = 3.48
$ws.Cells.Item($newRow, 15).Value = "04/11/2023 03:43"
$ws.Cells.Item($newRow, 16).Value = 3.14
$ws.Cells.Item($newRow, 17).Value = "05/11/2023 15:17"
$ws.Cells.Item($newRow, 18).Value = 5.17
$ws.Cells.Item($newRow, 19).Value = "04/11/2023 03:43"
$ws.Cells.Item($newRow, 20).Value = 4.41
$ws.Cells.Item($newRow, 21).Value = "05/11/2023 15:17"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/sabah-baku-zira-fk/QBc1Wk6o/"
